$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 keeps the same text content ("Run #1: METM - 8mths ALL ENABLED...") -
# no value changes needed there; shared-string renumbering happens automatically
# once the two obsolete strings below are removed from use.

# Row 13 previously held only a label ("Run #1: METM - 1 year NONE ENABLED...").
# Replace it with a brand new data row.
$ws.Range("A13").Value = "Run #1: METM - 9mths NONE ENABLED (20.06.20199)"
$ws.Range("B13").Value = 308000
$ws.Range("C13").Value = 0.03
$ws.Range("E13").Value = 6227
$ws.Range("F13").Value = 11737
$ws.Range("G13").Value = 3100
$ws.Range("H13").Value = 7690
$ws.Range("I13").Value = 7000
$ws.Range("J13").Value = 0.89
$ws.Range("K13").Value = "~350,000"
$ws.Range("L13").Value = 12200
$ws.Range("M13").Value = 0.40235919999999997
$ws.Range("N13").Value = "-     NEW TRANSFER MANAGER: NUM INVOCATIONS: 160056, TOTAL MS: 64400, AVG TIME/INVOCATION: 0.4023592ms"

# Match formatting of the row above (row 12) for the corresponding columns
$ws.Range("B13").NumberFormat = "#,##0"
$ws.Range("B13").Font.Color = 255

$ws.Range("C13").NumberFormat = "0%"

$ws.Range("F13").Font.Color = 255
$ws.Range("G13").Font.Color = 255

$ws.Range("J13").NumberFormat = "0%"
$ws.Range("J13").Font.Color = 255

$ws.Range("K13").Font.Color = 255

# Row 14 previously held only a label ("Run #1: METM - 1 year ONLY WH1st ENABLED...").
# It now holds a new label.
$ws.Range("A14").Value = "Run #2: VANILLA - 1 year (19.09.2199)"

$ws.Range("B14").Select() | Out-Null
